$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the remaining (first) data row with the corrected date/value pair
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 6.681483765882756

# Remove the now-obsolete trailing rows (3-17) entirely
$ws.Range("A3:A17").EntireRow.Delete()
